$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.574.53"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.11%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.434.91"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.10%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "567.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.97%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.59"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.44%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.05%  "

# Row 8
$ws.Range("E8").Value = "  +0.44%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.111"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.86%  "

# Row 10
$ws.Range("E10").Value = "  +0.44%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.31"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.38%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.355"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.92%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "26.87"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.47%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000180"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.78%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.860.63"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.69%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "62.438.39"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.05%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.428.91"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.33%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.24"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.21%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.98"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.39%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "324.05"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.04%  "

# Row 21
$ws.Range("E21").Value = "  +1.29%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.26%  "

# Row 23
$ws.Range("B23").Value = "SuiNetwork"
$ws.Range("C23").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.83"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.68%  "

# Row 24
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "67.25"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.69%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "593.58"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.40%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.57"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.21%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0₃0999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +7.26%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.550.19"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.38%  "

# Row 29
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.46"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.62%  "

# Row 30
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.15%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.45"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.54%  "

# Row 32
$ws.Range("E32").Value = "  -0.78%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.88"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.79%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.51"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.16%  "

# Row 35
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.85"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.99%  "

# Row 36
$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.05%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.382"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.76%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.75"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.46%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.36"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.51%  "

# Row 40
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.83"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.39%  "

# Row 41
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "147.83"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.94%  "

# Row 42
$ws.Range("E42").Value = "  +0.09%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.45"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +9.42%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "148.62"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.46%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.67"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.28%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0536"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.60%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "20.56"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.67%  "

# Row 48
$ws.Range("E48").Value = "  +2.00%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0231"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.02%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0919"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.31%  "

# Row 51
$ws.Range("E51").Value = "  +4.55%  "
